# "Corporate Marketing Line Fix in P202"
#
# The header in B1 used to read "Brand Category" (a string that was only
# ever used in that one header cell). It gets replaced with a new header,
# "Corporate Marketing Line" - everything else on the sheet (the
# &=result.* data-row formulas, all the other headers, styles, etc.) is
# left exactly as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the B1 header. Assigning the new text automatically drops the
# now-unused "Brand Category" shared string and appends the new one.
$ws.Range("B1").Value = "Corporate Marketing Line"

# Column B needs to be a bit wider to comfortably fit the new, longer
# header text (19.375 -> 25.125 characters).
$ws.Columns.Item(2).ColumnWidth = 24.36

# Leave the cursor where the editor left it after making the change.
$ws.Range("E11").Select()
